# Add "write data to general log file" (H1/H2 summary cell + two new eq_log rows)
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("main")

# --- New column H: "середне значення циклів між замінами" header + value ---
$ws.Columns.Item(8).ColumnWidth = 14.45

$ws.Range("H1").Value = "середне значення циклів між замінами"
$ws.Range("G1").Copy()
$ws.Range("H1").PasteSpecial(-4122)
$ws.Range("H1").Value = "середне значення циклів між замінами"

$ws.Range("H2").Value = 150000
$ws.Range("H2").NumberFormat = "#,##0"
$ws.Range("H2").HorizontalAlignment = -4108
$ws.Range("H2").VerticalAlignment = -4108
$ws.Range("H2").WrapText = $true

# --- New log entries: push the "**" end marker down two rows and fill in
#     the two new equipment-log rows above it ---
$marker = $ws.Range("A45").Value2

$ws.Range("A45").Value = "'08/02/2018"
$ws.Range("B45").Value = "'3012"
$ws.Range("C45").Value = "Рекваліфікація / EMPB"
$ws.Range("D45").Value = "'1335"
$ws.Range("A45:D45").Style = "Обычный"

$ws.Range("A46").Value = "'08/02/2018"
$ws.Range("B46").Value = "'3012"
$ws.Range("C46").Value = "Пошкодження поверхні контакту"
$ws.Range("D46").Value = "'1336"
$ws.Range("A46:D46").Style = "Обычный"

$ws.Range("A47").Value = $marker
$ws.Range("A47").Style = "Обычный"

# --- Selection / scroll position used when the file was saved ---
$ws.Activate()
$ws.Range("H1:H2").Select()
$excel.ActiveWindow.ScrollRow = 1
$excel.ActiveWindow.ScrollColumn = 1
